$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object "object[,]" 24,5
$arrBF[0,0] = 3.975606546170866
$arrBF[0,1] = 0.2864795845894434
$arrBF[0,2] = 0.008536791026919843
$arrBF[0,3] = 0.04218811757395358
$arrBF[0,4] = 4.778860952708982
$arrBF[1,0] = 3.881884048324764
$arrBF[1,1] = 0.2634821436224399
$arrBF[1,2] = 0.007445546137937242
$arrBF[1,3] = 0.04177184720871985
$arrBF[1,4] = 4.75714298633109
$arrBF[2,0] = 3.826681723177614
$arrBF[2,1] = 0.2495364413016432
$arrBF[2,2] = 0.006775747241810848
$arrBF[2,3] = 0.04151075920715286
$arrBF[2,4] = 4.746015564320786
$arrBF[3,0] = 3.804774224173684
$arrBF[3,1] = 0.2438968920009472
$arrBF[3,2] = 0.006502784756982294
$arrBF[3,3] = 0.0414029662815425
$arrBF[3,4] = 4.742034781527167
$arrBF[4,0] = 3.801171967153721
$arrBF[4,1] = 0.242963057561866
$arrBF[4,2] = 0.006457456608174539
$arrBF[4,3] = 0.04138498249888745
$arrBF[4,4] = 4.741407180372377
$arrBF[5,0] = 3.826383892514514
$arrBF[5,1] = 0.2494602091695697
$arrBF[5,2] = 0.006772066124682397
$arrBF[5,3] = 0.04150931115208678
$arrBF[5,4] = 4.745959637885818
$arrBF[6,0] = 3.942803682930844
$arrBF[6,1] = 0.278513459491819
$arrBF[6,2] = 0.008160438471179532
$arrBF[6,3] = 0.04204571897044218
$arrBF[6,4] = 4.770913588035413
$arrBF[7,0] = 4.189782338302166
$arrBF[7,1] = 0.3369008100935673
$arrBF[7,2] = 0.0108882102009531
$arrBF[7,3] = 0.0430547270735735
$arrBF[7,4] = 4.837433824306856
$arrBF[8,0] = 4.382768886533199
$arrBF[8,1] = 0.3807021077585659
$arrBF[8,2] = 0.01290027029359919
$arrBF[8,3] = 0.04377095092725902
$arrBF[8,4] = 4.897136646087603
$arrBF[9,0] = 4.473100869906773
$arrBF[9,1] = 0.4008342669047238
$arrBF[9,2] = 0.01381839404301388
$arrBF[9,3] = 0.04409156101893696
$arrBF[9,4] = 4.926673709018758
$arrBF[10,0] = 4.507674978684861
$arrBF[10,1] = 0.4084882060025166
$arrBF[10,2] = 0.01416655371483699
$arrBF[10,3] = 0.04421223720679723
$arrBF[10,4] = 4.93820244210707
$arrBF[11,0] = 4.500212457105533
$arrBF[11,1] = 0.4068384347849019
$arrBF[11,2] = 0.01409154840222726
$arrBF[11,3] = 0.04418627971910016
$arrBF[11,4] = 4.935704209360807
$arrBF[12,0] = 4.475937929409099
$arrBF[12,1] = 0.4014633497986893
$arrBF[12,2] = 0.01384702714642572
$arrBF[12,3] = 0.04410150369536936
$arrBF[12,4] = 4.927615286687626
$arrBF[13,0] = 4.461116976903554
$arrBF[13,1] = 0.3981749225540625
$arrBF[13,2] = 0.01369731650147799
$arrBF[13,3] = 0.04404948110785423
$arrBF[13,4] = 4.922705398717284
$arrBF[14,0] = 4.376916779492888
$arrBF[14,1] = 0.3793906347693223
$arrBF[14,2] = 0.01284033242303906
$arrBF[14,3] = 0.04374989498914061
$arrBF[14,4] = 4.895254328031285
$arrBF[15,0] = 4.325914935093806
$arrBF[15,1] = 0.3679204543806236
$arrBF[15,2] = 0.0123153809605796
$arrBF[15,3] = 0.04356478742373859
$arrBF[15,4] = 4.879024262666405
$arrBF[16,0] = 4.296819118230985
$arrBF[16,1] = 0.3613425317853967
$arrBF[16,2] = 0.01201370409262381
$arrBF[16,3] = 0.04345782765806838
$arrBF[16,4] = 4.869912875786099
$arrBF[17,0] = 4.287008782698194
$arrBF[17,1] = 0.3591186778293149
$arrBF[17,2] = 0.01191160396495405
$arrBF[17,3] = 0.04342152812888589
$arrBF[17,4] = 4.866866287940468
$arrBF[18,0] = 4.331319409556272
$arrBF[18,1] = 0.3691394602705032
$arrBF[18,2] = 0.01237123531851836
$arrBF[18,3] = 0.04358454311671345
$arrBF[18,4] = 4.880728814435287
$arrBF[19,0] = 4.483057961540794
$arrBF[19,1] = 0.4030413146417118
$arrBF[19,2] = 0.01391883513507963
$arrBF[19,3] = 0.04412642420052393
$arrBF[19,4] = 4.929981857672232
$arrBF[20,0] = 4.584369978590303
$arrBF[20,1] = 0.425375193608204
$arrBF[20,2] = 0.01493315219597235
$arrBF[20,3] = 0.04447631649002126
$arrBF[20,4] = 4.964175683777171
$arrBF[21,0] = 4.530101204993741
$arrBF[21,1] = 0.413438780226727
$arrBF[21,2] = 0.01439150286925184
$arrBF[21,3] = 0.04428995656826817
$arrBF[21,4] = 4.94574183283612
$arrBF[22,0] = 4.328875343405912
$arrBF[22,1] = 0.3685882965559699
$arrBF[22,2] = 0.01234598318099955
$arrBF[22,3] = 0.04357561324790993
$arrBF[22,4] = 4.879957503037161
$arrBF[23,0] = 4.120953628334291
$arrBF[23,1] = 0.3209497907135983
$arrBF[23,2] = 0.01014922713522282
$arrBF[23,3] = 0.04278624686196775
$arrBF[23,4] = 4.81754446409326
$ws.Range("B2:F25").Value = $arrBF

$arrJ = New-Object "object[,]" 24,1
$arrJ[0,0] = 0.1229205903000947
$arrJ[1,0] = 0.1231051025490411
$arrJ[2,0] = 0.1232238719411609
$arrJ[3,0] = 0.1232736517025748
$arrJ[4,0] = 0.1232820010404563
$arrJ[5,0] = 0.1232245376957071
$arrJ[6,0] = 0.1229830748070286
$arrJ[7,0] = 0.1225529035411514
$arrJ[8,0] = 0.1222630884284066
$arrJ[9,0] = 0.1221369000660992
$arrJ[10,0] = 0.1220899254596408
$arrJ[11,0] = 0.1221000062835795
$arrJ[12,0] = 0.1221330192121779
$arrJ[13,0] = 0.1221533460279627
$arrJ[14,0] = 0.1222714486312926
$arrJ[15,0] = 0.1223453462207478
$arrJ[16,0] = 0.1223883820968061
$arrJ[17,0] = 0.1224030447355902
$arrJ[18,0] = 0.1223374246593067
$arrJ[19,0] = 0.1221233005444011
$arrJ[20,0] = 0.1219880794196513
$arrJ[21,0] = 0.1220598182156607
$arrJ[22,0] = 0.1223410042804325
$arrJ[23,0] = 0.1226646562850209
$ws.Range("J2:J25").Value = $arrJ

$arrLN = New-Object "object[,]" 24,3
$arrLN[0,0] = 0.3265923476159784
$arrLN[0,1] = 0.7321708664111171
$arrLN[0,2] = 3.259197752174316
$arrLN[1,0] = 0.3253257406725965
$arrLN[1,1] = 0.7192830793392631
$arrLN[1,2] = 3.27552583936469
$arrLN[2,0] = 0.3246897639083883
$arrLN[2,1] = 0.7117789697486359
$arrLN[2,2] = 3.286377293793137
$arrLN[3,0] = 0.3244662581945406
$arrLN[3,1] = 0.7088236774714929
$arrLN[3,2] = 3.29100657917644
$arrLN[4,0] = 0.3244312997037611
$arrLN[4,1] = 0.7083391512283441
$arrLN[4,2] = 3.291787771275409
$arrLN[5,0] = 0.3246866052155113
$arrLN[5,1] = 0.7117386980587384
$arrLN[5,2] = 3.286438887534331
$arrLN[6,0] = 0.3261262180837434
$arrLN[6,1] = 0.7276421536213817
$arrLN[6,2] = 3.264655974072667
$arrLN[7,0] = 0.3300734479392986
$arrLN[7,1] = 0.7620846647675421
$arrLN[7,2] = 3.228514076392102
$arrLN[8,0] = 0.3336591854330351
$arrLN[8,1] = 0.7893929252097962
$arrLN[8,2] = 3.205997497605551
$arrLN[9,0] = 0.335439502659213
$arrLN[9,1] = 0.8022555433902738
$arrLN[9,2] = 3.196637113978909
$arrLN[10,0] = 0.3361351122405694
$arrLN[10,1] = 0.8071898374477016
$arrLN[10,2] = 3.193220043669754
$arrLN[11,0] = 0.3359843466867858
$arrLN[11,1] = 0.8061243209282125
$arrLN[11,2] = 3.19395029048313
$arrLN[12,0] = 0.3354963011681917
$arrLN[12,1] = 0.8026602167245613
$arrLN[12,2] = 3.196353430719441
$arrLN[13,0] = 0.3352001514740692
$arrLN[13,1] = 0.8005466295438097
$arrLN[13,2] = 3.197842046195817
$arrLN[14,0] = 0.3335458382541816
$arrLN[14,1] = 0.7885611957210088
$arrLN[14,2] = 3.20662703017986
$arrLN[15,0] = 0.3325691665526449
$arrLN[15,1] = 0.7813213671413664
$arrLN[15,2] = 3.212242790549539
$arrLN[16,0] = 0.3320214494425784
$arrLN[16,1] = 0.7771986051275377
$arrLN[16,2] = 3.215555858055154
$arrLN[17,0] = 0.3318384131767544
$arrLN[17,1] = 0.7758098115236081
$arrLN[17,2] = 3.216691848299092
$arrLN[18,0] = 0.3326716820224362
$arrLN[18,1] = 0.7820877732635552
$arrLN[18,2] = 3.211636385614227
$arrLN[19,0] = 0.3356390700653122
$arrLN[19,1] = 0.8036759828174596
$arrLN[19,2] = 3.195644104938907
$arrLN[20,0] = 0.3377034097486558
$arrLN[20,1] = 0.8181553663515473
$arrLN[20,2] = 3.18593570000948
$arrLN[21,0] = 0.3365901982488708
$arrLN[21,1] = 0.8103934897450955
$arrLN[21,2] = 3.191049029965455
$arrLN[22,0] = 0.3326252918427883
$arrLN[22,1] = 0.7817411580025464
$arrLN[22,2] = 3.211910278156822
$arrLN[23,0] = 0.3288852361706773
$arrLN[23,1] = 0.7524166107989032
$arrLN[23,2] = 3.237584629830593
$ws.Range("L2:N25").Value = $arrLN

